$d = $word.ActiveDocument

# The template paragraph holds a single Word field whose code is an M2Doc
# query (e.g. "m:self.asImageByRepresentationDescriptionName('Entities')
# ->first().setWidth(300)"). The parser was changed to work off plain
# "{...}" text tokens instead of real Word fields, so this field needs to
# be turned into literal text: the field characters/instrText runs are
# replaced by plain <w:t> runs holding the same text, wrapped in "{" and
# "}" so it reads as "{m:...}".

$f = $d.Fields.Item(1)
$code = $f.Code.Text

# Locate the paragraph that owns the field before it gets removed.
# (Field.Code.Paragraphs.Item(1).Index is relative to that sub-range, not
# the document, so walk Document.Paragraphs and match on position.)
$fieldPos = $f.Code.Start
$pIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if (($fieldPos -ge $candidate.Range.Start) -and ($fieldPos -lt $candidate.Range.End)) {
        $pIndex = $i
        break
    }
}

# Split the field code exactly the way the original template's runs did,
# so every fragment becomes its own run (matching the target markup).
$parts = New-Object System.Collections.ArrayList
[void]$parts.Add("{m:")
[void]$parts.Add("self")
[void]$parts.Add(".")
[void]$parts.Add("asImageByRepresentationDescriptionName")
[void]$parts.Add("(")
[void]$parts.Add("'")
[void]$parts.Add("Entities")
[void]$parts.Add("'")
[void]$parts.Add(")")
[void]$parts.Add("->first()")
[void]$parts.Add(".setWidth(300)}")

function Escape-Xml([string]$s) {
    $s = $s.Replace("&", "&amp;")
    $s = $s.Replace("<", "&lt;")
    $s = $s.Replace(">", "&gt;")
    return $s
}

$runsXml = ""
foreach ($part in $parts) {
    $runsXml += "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>" + (Escape-Xml $part) + "</w:t></w:r>"
}

# Remove the field (begin/instrText.../end) entirely, leaving the
# paragraph mark behind.
$f.Delete()

$p = $d.Paragraphs.Item($pIndex)
$rng = $p.Range
$rng.Collapse(1)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidP="00340939" w:rsidR="00A10D75" w:rsidRDefault="00474E78"><w:pPr><w:widowControl w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($xml)
